# Added session management constraints for all routes
$wb = $excel.ActiveWorkbook

# --- Hackathon sheet: update winning_team_profiles (jury) string for first hackathon row ---
$wsHack = $wb.Worksheets.Item("Hackathon")
$wsHack.Range("M2").Value = "[['Mr. Aditya', 'ISL'], ['Mr. Manoj', 'GBS'], ['Dr. Dilip', 'GTS']]"

# --- Tech Session sheet: update existing data row with new session details ---
$wsTech = $wb.Worksheets.Item("Tech Session")
$wsTech.Range("C2").Value = "tech sess 1"
$wsTech.Range("D2").Value = "IISc"

# no_of_participants is stored as text in this workbook, force text so the
# numeric-looking value doesn't get auto-converted to a number
$wsTech.Range("E2").NumberFormat = "@"
$wsTech.Range("E2").Value = "200000"
$wsTech.Range("E2").Style = "Normal"

$wsTech.Range("F2").Value = "Planned"

# startdate / enddate are stored as text dates, force text so they are not
# auto-converted to date serial numbers
$wsTech.Range("G2").NumberFormat = "@"
$wsTech.Range("G2").Value = "2020-10-31"
$wsTech.Range("G2").Style = "Normal"

$wsTech.Range("H2").NumberFormat = "@"
$wsTech.Range("H2").Value = "2020-10-31"
$wsTech.Range("H2").Style = "Normal"

$wsTech.Range("J2").Value = "ISL"
$wsTech.Range("K2").Value = "SME_NAME_2"

# --- SUR sheet: add new data row (row 2) with session management constraints ---
$wsSur = $wb.Worksheets.Item("SUR")

# Match header row formatting (bold, centered, bordered) for the index cell in column A
$wsTech.Range("A2").Copy($wsSur.Range("A2"))
$wsSur.Range("A2").Value = 0

$wsSur.Range("B2").Value = "ABC"
$wsSur.Range("C2").Value = "ABC"
$wsSur.Range("D2").Value = "tech #1"

$wsSur.Range("E2").NumberFormat = "@"
$wsSur.Range("E2").Value = "2020-11-01"
$wsSur.Range("E2").Style = "Normal"

$wsSur.Range("F2").NumberFormat = "@"
$wsSur.Range("F2").Value = "2020-11-01"
$wsSur.Range("F2").Style = "Normal"

$wsSur.Range("G2").NumberFormat = "@"
$wsSur.Range("G2").Value = "2020-11-01"
$wsSur.Range("G2").Style = "Normal"

$wsSur.Range("H2").NumberFormat = "@"
$wsSur.Range("H2").Value = "2020-11-01"
$wsSur.Range("H2").Style = "Normal"

$wsSur.Range("I2").Value = "Approved"

$wsSur.Range("J2").NumberFormat = "@"
$wsSur.Range("J2").Value = "2020-11-01"
$wsSur.Range("J2").Style = "Normal"

$wsSur.Range("K2").NumberFormat = "@"
$wsSur.Range("K2").Value = "2020-11-01"
$wsSur.Range("K2").Style = "Normal"

$wsSur.Range("L2").Value = "Springer"
$wsSur.Range("M2").Value = "Yes"
$wsSur.Range("N2").Value = "SUR CITY"
$wsSur.Range("O2").Value = "aSCHAJK"
